{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block, along\n// with the blank paragraph that precedes it, from the Bibliografia section.\n// Target text that anchors the deletion:\n//   - (blank paragraph right after the \"...Bertero...\" reference)\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//      pages. Original theme under Creative Commons Attribution\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter...\" paragraph by its text content; this is the\n// stable anchor for the block that needs to be removed.\nconst jupiterIndex = items.findIndex((p) =>\n  p.text.indexOf(\"Ver no Jupiter\") !== -1\n);\n\nif (jupiterIndex === -1) {\n  throw new Error(\"Could not find the 'Ver no Jupiter...' paragraph.\");\n}\n\n// The paragraph right before it should be the blank spacer paragraph, and\n// the paragraph right after it should be the copyright/footer paragraph.\nconst toDelete = [];\nif (jupiterIndex - 1 >= 0 && items[jupiterIndex - 1].text === \"\") {\n  toDelete.push(items[jupiterIndex - 1]);\n}\ntoDelete.push(items[jupiterIndex]);\nif (\n  jupiterIndex + 1 < items.length &&\n  items[jupiterIndex + 1].text.indexOf(\"Powered by Jekyll\") !== -1\n) {\n  toDelete.push(items[jupiterIndex + 1]);\n}\n\n// Delete from last to first so indices/anchors stay valid as we go.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block, along\n# with the blank paragraph that precedes it, from the Bibliografia section.\n# Target text that anchors the deletion:\n#   - (blank paragraph right after the \"...Bertero...\" reference)\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#      pages. Original theme under Creative Commons Attribution\"\n\n$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter...\" paragraph by its text content; this is the\n# stable anchor for the block that needs to be removed.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"Ver no Jupiter\")\nif (-not $found) {\n    throw \"Could not find the 'Ver no Jupiter...' paragraph.\"\n}\n$searchRange.Expand(4) | Out-Null   # wdParagraph: grow the hit to its full paragraph\n\n# Resolve which paragraph (by collection index) that is, so we can inspect\n# its neighbours.\n$allParas = $d.Paragraphs\n$jupiterIndex = -1\nfor ($i = 1; $i -le $allParas.Count; $i++) {\n    if ($allParas.Item($i).Range.Start -eq $searchRange.Start) {\n        $jupiterIndex = $i\n        break\n    }\n}\nif ($jupiterIndex -eq -1) {\n    throw \"Could not resolve the paragraph index for the 'Ver no Jupiter...' paragraph.\"\n}\n\n$startPara = $allParas.Item($jupiterIndex)\n$endPara = $allParas.Item($jupiterIndex)\n\n# Include the blank spacer paragraph immediately before it, if present.\nif ($jupiterIndex -gt 1) {\n    $prevPara = $allParas.Item($jupiterIndex - 1)\n    if ($prevPara.Range.Text.Trim().Length -eq 0) {\n        $startPara = $prevPara\n    }\n}\n\n# Include the copyright/footer paragraph immediately after it, if present.\nif ($jupiterIndex -lt $allParas.Count) {\n    $nextPara = $allParas.Item($jupiterIndex + 1)\n    if ($nextPara.Range.Text -like \"*Powered by Jekyll*\") {\n        $endPara = $nextPara\n    }\n}\n\n$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$delRange.Delete()\n"}
